$wb = $excel.ActiveWorkbook

# --- Sheet "Metadata": update Date value ---
$wsMeta = $wb.Worksheets.Item("Metadata")
$wsMeta.Range("B8").Value = "2025-07-21T11:52:46+00:00"

# --- Sheet "Elements": update Binding Value Set URLs ---
$wsElem = $wb.Worksheets.Item("Elements")
$wsElem.Range("Z3").Value = "https://interop.esante.gouv.fr/terminologies/CodeSystem-TRE-R14-TypeDiplome?vs"
$wsElem.Range("Z4").Value = "https://interop.esante.gouv.fr/terminologies/CodeSystem-TRE-R16-LieuFormation?vs"

# --- Sheet "Elements": widen column Z (26) to fit new content ---
# (target OOXML width is 67.49609375; the COM width grid only allows
# 1/6-character increments, so 66.6665 is the input that rounds to the
# closest achievable stored width, 67.5)
$wsElem.Columns.Item(26).ColumnWidth = 66.6665
